$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: new diary entry (01.12.18 - Open External Software / Tooltipps) ---

# A33 must end up as a shared-string "01.12.18" while KEEPING the existing
# date-column style (numFmtId 164, no quotePrefix) - exactly like the other
# date cells in column A (e.g. A32 "28.11.18"). Typing the text straight into
# the cell makes Excel treat it as a real date serial and mint a new style
# with quotePrefix, so instead: copy an existing text-date cell's
# format+type down into A33 (preserves style "3" / t="s"), then overwrite
# the *value only* by pasting in literal text staged on a scratch cell.
$ws.Range("H1").Value = "'01.12.18"
$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial(-4163)   # xlPasteValues -> copies type+text, keeps A33's own style
$ws.Range("H1").Copy()
$ws.Range("A33").PasteSpecial(-4163)   # xlPasteValues -> swap in the real "01.12.18" text
$ws.Range("H1").Clear()                # wipe the scratch cell completely (no residue)

# B33 / C33: begin / end times
$ws.Range("B33").Value = 0.47916666666666669
$ws.Range("C33").Value = 0.79166666666666663

# D33: elapsed time formula (same shape as the rest of column D)
$ws.Range("D33").Formula = "=C33-B33"

# E33: remarks, wrapped two-line text - copy the wrap-text style used by the
# other multi-line remark cells (e.g. E9) so E33 lands on the same style
# index as the target, then fill in the real text.
$ws.Range("E9").Copy()
$ws.Range("E33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E33").Value = "-Open External Software`n-Tooltipps"

# Row height: two wrapped lines -> 30pt (matches the other 2-line remark rows)
$ws.Rows("33").RowHeight = 30

# --- Selection moved to E34 ---
$ws.Range("E34").Select()

$wb.Save()
